$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row (row 1) to short machine-friendly column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Normalize Spanish place-name casing (capitalize connector words: De, Del, El, La, Las, Los, Y)
#    and fix two special-case names (GUANAJUATO -> Guanajuato, MonteMorelos -> Montemorelos)
$textUpdates = @{
    'B4' = 'Pabellón De Arteaga'
    'B5' = 'Rincón De Romos'
    'B22' = 'Benemérito De Las Américas'
    'B65' = 'Guadalupe Y Calvo'
    'B66' = 'Hidalgo Del Parral'
    'B80' = 'San Juan De Sabinas'
    'B86' = 'Villa De Álvarez'
    'A88' = 'Ciudad De México'
    'B92' = 'Cuajimalpa De Morelos'
    'B112' = 'Nombre De Dios'
    'B113' = 'Pánuco De Coronado'
    'A120' = 'Estado De México'
    'B120' = 'Acambay De Ruíz Castañeda'
    'B121' = 'Almoloya De Juárez'
    'B125' = 'Coacalco De Berriozábal'
    'B127' = 'Ecatepec De Morelos'
    'B133' = 'Naucalpan De Juárez'
    'B138' = 'San Simón De Guerero'
    'B143' = 'Tenango Del Valle'
    'B144' = 'Tlalnepantla De Baz'
    'B147' = 'Villa De Allende'
    'B148' = 'Villa Del Carbón'
    'A151' = 'Guanajuato'
    'B154' = 'Apaseo El Alto'
    'B155' = 'Apaseo El Grande'
    'B161' = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    'B164' = 'Jaral Del Progreso'
    'B173' = 'San Diego De La Unión'
    'B176' = 'San Luis De La Paz'
    'B177' = 'Santa Cruz De Juventino Rosas'
    'B180' = 'Valle De Santiago'
    'B186' = 'Acapulco De Juárez'
    'B188' = 'Ajuchitlán Del Progreso'
    'B189' = 'Alcozauca De Guerero'
    'B192' = 'Atenango Del Río'
    'B193' = 'Ayutla De Los Libres'
    'B195' = 'Buenavista De Cuéllar'
    'B196' = 'Chilapa De Álvarez'
    'B197' = 'Chilpancingo De Los Bravo'
    'B200' = 'Coyuca De Benítez'
    'B201' = 'Coyuca De Catalán'
    'B203' = 'Cuetzala Del Progreso'
    'B204' = 'Cutzamala De Pinzón'
    'B209' = 'Huitzuco De Los Figueroa'
    'B210' = 'Zihuatanejo De Azueta'
    'B217' = 'Taxco De Alarcón'
    'B218' = 'Técpan De Galeana'
    'B220' = 'Tepecoacuilco De Trujano'
    'B221' = 'Tixtla De Guerero'
    'B223' = 'Tlapa De Comonfort'
    'B228' = 'Atotonilco El Grande'
    'B235' = 'Jacala De Ledezma'
    'B238' = 'Mixquiahuala De Juárez'
    'B239' = 'Pachuca De Soto'
    'B246' = 'Tenango De Doria'
    'B247' = 'Tepehuacán De Guerero'
    'B248' = 'Tepeji Del Río De Ocampo'
    'B249' = 'Tezontepec De Aldama'
    'B254' = 'Tulancingo De Bravo'
    'B255' = 'Zacualtipán De Ángeles'
    'B262' = 'Autlán De Navarro'
    'B265' = 'Encarnación De Díaz'
    'B270' = 'Lagos De Moreno'
    'B275' = 'San Miguel El Alto'
    'B276' = 'Talpa De Allende'
    'B277' = 'Tamazula De Gordiano'
    'B281' = 'Tepatitlán De Morelos'
    'B282' = 'Tizapán El Alto'
    'B287' = 'Zapotitlán De Vadillo'
    'B318' = 'Tiquicheo De Nicolás Romero'
    'B336' = 'Tetela Del Volcán'
    'B337' = 'Tlaltizapán De Zapata'
    'B345' = 'Santa María Del Oro'
    'B357' = 'Montemorelos'
    'B359' = 'San Nicolás De Los Garza'
    'B362' = 'Acatlán De Pérez Figueroa'
    'B364' = 'Constancia Del Rosario'
    'B365' = 'Heroica Ciudad De Ejutla De Crespo'
    'B366' = 'Heroica Ciudad De Huajuapan De León'
    'B367' = 'Heroica Ciudad De Juchitán De Zaragoza'
    'B369' = 'Mazatlán Villa De Flores'
    'B370' = 'Miahuatlán De Porfirio Díaz'
    'B372' = 'Oaxaca De Juárez'
    'B373' = 'Ocotlán De Morelos'
    'B400' = 'Santo Domingo De Morelos'
    'B403' = 'Teotitlán De Flores Magón'
    'B404' = 'Totontepec Villa De Morelos'
    'B405' = 'Villa De Etla'
    'B406' = 'Villa De Tututepec De Melchor Ocampo'
    'B407' = 'Zimatlán De Álvarez'
    'B419' = 'Cuetzalan Del Progreso'
    'B423' = 'Huehuetlán El Chico'
    'B431' = 'Palmar De Bravo'
    'B437' = 'San Salvador El Seco'
    'B439' = 'Tepatlaxco De Hidalgo'
    'B441' = 'Tepexi De Rodríguez'
    'B442' = 'Tlacotepec De Benito Juárez'
    'B444' = 'Tuzamapan De Galeana'
    'B455' = 'Amealco De Bonfil'
    'B457' = 'Cadereyta De Montes'
    'B462' = 'Jalpan De Serra'
    'B463' = 'Landa De Matamoros'
    'B465' = 'Pinal De Amoles'
    'B467' = 'San Juan Del Río'
    'B479' = 'Cerro De San Pedro'
    'B481' = 'Ciudad Del Maíz'
    'B490' = 'Mexquitic De Carmona'
    'B503' = 'Tanquián De Escobedo'
    'B505' = 'Villa De Arista'
    'B506' = 'Villa De Arriaga'
    'B507' = 'Villa De Guadalupe'
    'B508' = 'Villa De Ramos'
    'B509' = 'Villa De Reyes'
    'B558' = 'Soto La Marina'
    'B577' = 'Boca Del Río'
    'B586' = 'Cosamaloapan De Carpio'
    'B593' = 'Hueyapan De Ocampo'
    'B594' = 'Ignacio De La Llave'
    'B596' = 'Ixhuatlán Del Café'
    'B600' = 'Juchique De Ferrer'
    'B605' = 'Martínez De La Torre'
    'B606' = 'Medellín De Bravo'
    'B614' = 'Ozuluama De Mascareñas'
    'B617' = 'Paso De Ovejas'
    'B620' = 'Poza Rica De Hidalgo'
    'B627' = 'Soledad De Doblado'
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# 3) Correct the floating point precision of percentage values (re-derived computation)
$numUpdates = @{
    'D4' = 0.0009708737864077668
    'D9' = 0.0009708737864077668
    'D15' = 0.0009708737864077668
    'D20' = 0.0009708737864077668
    'D22' = 0.0009708737864077668
    'D23' = 0.0009708737864077668
    'D27' = 0.0009708737864077668
    'D28' = 0.0009708737864077668
    'D29' = 0.0009708737864077668
    'D31' = 0.0009708737864077668
    'D33' = 0.0009708737864077668
    'D35' = 0.0009708737864077668
    'D38' = 0.0009708737864077668
    'D44' = 0.0009708737864077668
    'D53' = 0.0009708737864077668
    'D61' = 0.0009708737864077668
    'D64' = 0.0009708737864077668
    'D69' = 0.0009708737864077668
    'D75' = 0.0009708737864077668
    'D87' = 0.0009708737864077668
    'D90' = 0.0009708737864077668
    'D92' = 0.0009708737864077668
    'D113' = 0.0009708737864077668
    'D114' = 0.0009708737864077668
    'D117' = 0.0009708737864077668
    'D133' = 0.0009708737864077668
    'D147' = 0.0009708737864077668
    'D152' = 0.0009708737864077668
    'D158' = 0.0009708737864077668
    'D160' = 0.0009708737864077668
    'D164' = 0.0009708737864077668
    'D167' = 0.0009708737864077668
    'D169' = 0.0009708737864077668
    'D170' = 0.0009708737864077668
    'D176' = 0.009708737864077667
    'D194' = 0.0009708737864077668
    'D204' = 0.0009708737864077668
    'D211' = 0.0009708737864077668
    'D212' = 0.0009708737864077668
    'D214' = 0.0009708737864077668
    'D218' = 0.0009708737864077668
    'D231' = 0.0009708737864077668
    'D239' = 0.0009708737864077668
    'D245' = 0.0009708737864077668
    'D249' = 0.0009708737864077668
    'D250' = 0.0009708737864077668
    'D260' = 0.0009708737864077668
    'D270' = 0.0009708737864077668
    'D272' = 0.0009708737864077668
    'D273' = 0.0009708737864077668
    'D278' = 0.0009708737864077668
    'D281' = 0.0009708737864077668
    'D296' = 0.0009708737864077668
    'D300' = 0.0009708737864077668
    'D303' = 0.0009708737864077668
    'D306' = 0.0009708737864077668
    'D308' = 0.0009708737864077668
    'D315' = 0.0009708737864077668
    'D316' = 0.0009708737864077668
    'D320' = 0.0009708737864077668
    'D323' = 0.0009708737864077668
    'D329' = 0.0009708737864077668
    'D338' = 0.0009708737864077668
    'D340' = 0.0009708737864077668
    'D346' = 0.0009708737864077668
    'D347' = 0.0009708737864077668
    'D350' = 0.0009708737864077668
    'D351' = 0.0009708737864077668
    'D360' = 0.0009708737864077668
    'D366' = 0.0009708737864077668
    'D377' = 0.0009708737864077668
    'D379' = 0.0009708737864077668
    'D382' = 0.0009708737864077668
    'D385' = 0.0009708737864077668
    'D387' = 0.0009708737864077668
    'D393' = 0.0009708737864077668
    'D395' = 0.0009708737864077668
    'D396' = 0.0009708737864077668
    'D397' = 0.0009708737864077668
    'D398' = 0.0009708737864077668
    'D406' = 0.0009708737864077668
    'D411' = 0.0009708737864077668
    'D412' = 0.0009708737864077668
    'D422' = 0.0009708737864077668
    'D423' = 0.0009708737864077668
    'D426' = 0.0009708737864077668
    'D432' = 0.0009708737864077668
    'D438' = 0.0009708737864077668
    'D444' = 0.0009708737864077668
    'D446' = 0.0009708737864077668
    'D448' = 0.0009708737864077668
    'D451' = 0.0009708737864077668
    'D453' = 0.0009708737864077668
    'D459' = 0.0009708737864077668
    'D463' = 0.0009708737864077668
    'D465' = 0.0009708737864077668
    'D470' = 0.0009708737864077668
    'D474' = 0.0009708737864077668
    'D476' = 0.0009708737864077668
    'D477' = 0.0009708737864077668
    'D487' = 0.0009708737864077668
    'D495' = 0.0009708737864077668
    'D496' = 0.0009708737864077668
    'D501' = 0.0009708737864077668
    'D503' = 0.0009708737864077668
    'D516' = 0.0009708737864077668
    'D523' = 0.0009708737864077668
    'D534' = 0.0009708737864077668
    'D536' = 0.0009708737864077668
    'D537' = 0.009708737864077667
    'D545' = 0.0009708737864077668
    'D546' = 0.0009708737864077668
    'D554' = 0.0009708737864077668
    'D557' = 0.0009708737864077668
    'D568' = 0.0009708737864077668
    'D576' = 0.0009708737864077668
    'D579' = 0.0009708737864077668
    'D580' = 0.0009708737864077668
    'D583' = 0.0009708737864077668
    'D586' = 0.0009708737864077668
    'D589' = 0.0009708737864077668
    'D590' = 0.0009708737864077668
    'D591' = 0.0009708737864077668
    'D597' = 0.0009708737864077668
    'D602' = 0.0009708737864077668
    'D606' = 0.0009708737864077668
    'D610' = 0.0009708737864077668
    'D614' = 0.0009708737864077668
    'D622' = 0.0009708737864077668
    'D637' = 0.0009708737864077668
    'D638' = 0.0009708737864077668
    'D640' = 0.0009708737864077668
    'D646' = 0.0009708737864077668
    'D653' = 0.0009708737864077668
    'D656' = 0.0009708737864077668
    'D665' = 0.0009708737864077668
    'D667' = 0.0009708737864077668
}
foreach ($ref in $numUpdates.Keys) {
    $ws.Range($ref).Value = $numUpdates[$ref]
}

# 4) Remove trailing footnote/metadata rows (673-677) which are no longer part of the clean dataset
$ws.Range("A673:A677").EntireRow.Delete()

Write-Host "Done. UsedRange is now:" $ws.UsedRange.Address()